# Automatische test-sync: 2025-06-23 18:39:50
# Adds the new "Order wijzigen" mail-log entry (row 18) to the Logs sheet,
# extends the conditional formatting ranges accordingly, and updates the
# Dashboard summary counts (Bestelling / Levering now 3, Retour / Terugbetaling stays 3,
# with both rows swapping position to mirror the sort order produced by the sync).

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row ---------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A18").Value2 = "Order wijzigen"
$logs.Range("B18").Value2 = "mailmind.test@zohomail.eu"
$logs.Range("C18").Value2 = "Kan ik mijn bestelling nog aanpassen?"
$logs.Range("D18").Value2 = "Bestelling / Levering"
$logs.Range("F18").Value2 = "2025-06-23 18:39:29"
$logs.Range("G18").Value2 = "Nee"

# --- Extend conditional formatting ranges to cover the new row ------------
$dFormats = $logs.Range("D2:D17").FormatConditions
$dFormats.Item(1).ModifyAppliesToRange($logs.Range("D2:D18"))

$gFormats = $logs.Range("G2:G17").FormatConditions
$gFormats.Item(1).ModifyAppliesToRange($logs.Range("G2:G18"))

# --- Dashboard sheet: update category counts -------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value2 = "Bestelling / Levering"
$dash.Range("A4").Value2 = "Retour / Terugbetaling"
$dash.Range("B3").Value2 = 3
$dash.Range("B4").Value2 = 3
